$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 2-43) holds the "ID Competição" values which were
# incorrectly recorded as 50; they should be 250.
$ws.Range("B2:B43").Value = 250
